{"js": "// The revision captured by the diff only reorders/renames the XML\n// namespace-prefix declarations on the root element of document.xml,\n// endnotes.xml, footer1.xml, footer2.xml, footnotes.xml, header1.xml,\n// numbering.xml, styles.xml and theme1.xml (e.g. \"ns17\" -> \"ns19\", the\n// \"m\" prefix moving later in the attribute list, \"o\"/\"v\" swapping\n// order, ...). Every prefix still maps to the same namespace URI and\n// every other byte of every part (body text, tables, headers/footers,\n// styles, numbering, theme colors, etc.) is unchanged - this is the\n// kind of cosmetic churn produced when the file is re-saved by a\n// different Word build, not a reachable operation of the Word\n// JavaScript API (namespace-prefix serialization isn't something\n// Office.js exposes any control over).\n//\n// So there is no content to mutate here. Touch/read the body (without\n// changing it) so the sync round-trips cleanly and leaves the document\n// byte-for-byte equivalent in content.\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n", "ps1": "# The revision captured by the diff only reorders/renames the XML\n# namespace-prefix declarations on the root element of document.xml,\n# endnotes.xml, footer1.xml, footer2.xml, footnotes.xml, header1.xml,\n# numbering.xml, styles.xml and theme1.xml (e.g. \"ns17\" -> \"ns19\", the\n# \"m\" prefix moving later in the attribute list, \"o\"/\"v\" swapping\n# order, ...). Every prefix still maps to the same namespace URI and\n# every other byte of every part (body text, tables, headers/footers,\n# styles, numbering, theme colors, etc.) is unchanged - this is the\n# kind of cosmetic churn produced when the file is re-saved by a\n# different Word build, not something the Word COM object model has\n# any property/method for (namespace-prefix serialization isn't part\n# of the object model).\n#\n# So there is no content to mutate here. Touch/read the document\n# (without changing it) so the automation round-trips cleanly and\n# leaves the document byte-for-byte equivalent in content.\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
